# Auto-generated edit script applying the Goblin_Profits market-data refresh diff
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H40").Value = 2360.3635
$ws.Range("I40").Value = 1294
$ws.Range("J40").Value = 3640
$ws.Range("K40").Value = 1294
$ws.Range("L40").Value = 3640
$ws.Range("M40").Value = -1119
$ws.Range("N40").Value = -3990
$ws.Range("H61").Value = 787.25
$ws.Range("I61").Value = 787.25
$ws.Range("K61").Value = 2361.75
$ws.Range("M61").Value = -2189.75
$ws.Range("H64").Value = 8615.346
$ws.Range("I64").Value = 4400
$ws.Range("J64").Value = 9619
$ws.Range("K64").Value = 4400
$ws.Range("L64").Value = 9619
$ws.Range("M64").Value = -4152
$ws.Range("N64").Value = -10115
$ws.Range("H67").Value = 8615.346
$ws.Range("I67").Value = 4400
$ws.Range("J67").Value = 9619
$ws.Range("K67").Value = 4400
$ws.Range("L67").Value = 9619
$ws.Range("M67").Value = -3542
$ws.Range("N67").Value = -11335
$ws.Range("H107").Value = 1929.5
$ws.Range("J107").Value = 0
$ws.Range("L107").Value = 0
$ws.Range("N107").Value = $null
$ws.Range("H113").Value = 4214.2856
$ws.Range("I113").Value = 4000
$ws.Range("J113").Value = 5000
$ws.Range("K113").Value = 4000
$ws.Range("L113").Value = 5000
$ws.Range("M113").Value = -746
$ws.Range("N113").Value = -11508
$ws.Range("H132").Value = 1577.6945
$ws.Range("I132").Value = 934.28
$ws.Range("K132").Value = 2802.84
$ws.Range("M132").Value = -272.8400000000001
$ws.Range("H137").Value = 1433.56
$ws.Range("I137").Value = 1255.579
$ws.Range("J137").Value = 1997.1666
$ws.Range("K137").Value = 3766.737
$ws.Range("L137").Value = 5991.4998
$ws.Range("M137").Value = -1216.737
$ws.Range("N137").Value = -11091.4998
$ws.Range("H138").Value = 2959.1572
$ws.Range("I138").Value = 1985.7222
$ws.Range("J138").Value = 3296.1155
$ws.Range("K138").Value = 5957.1666
$ws.Range("L138").Value = 9888.3465
$ws.Range("M138").Value = -817.1665999999996
$ws.Range("N138").Value = -20168.3465

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H2").Value = 3128.4
$ws.Range("I2").Value = 99.5
$ws.Range("K2").Value = 99.5
$ws.Range("M2").Value = 13.5
$ws.Range("H32").Value = 4473.1
$ws.Range("J32").Value = 10720.5
$ws.Range("L32").Value = 10720.5
$ws.Range("N32").Value = -11294.5
$ws.Range("H45").Value = 1562.0714
$ws.Range("I45").Value = 1638.5454
$ws.Range("K45").Value = 1638.5454
$ws.Range("M45").Value = -1261.5454
$ws.Range("H63").Value = 7750
$ws.Range("J63").Value = 7500
$ws.Range("L63").Value = 7500
$ws.Range("N63").Value = -8872
$ws.Range("H66").Value = 7750
$ws.Range("J66").Value = 7500
$ws.Range("L66").Value = 37500
$ws.Range("N66").Value = -44364
$ws.Range("H74").Value = 1471.4584
$ws.Range("I74").Value = 1448.5217
$ws.Range("K74").Value = 1448.5217
$ws.Range("M74").Value = -574.5217
$ws.Range("H77").Value = 1471.4584
$ws.Range("I77").Value = 1448.5217
$ws.Range("K77").Value = 7242.6085
$ws.Range("M77").Value = -2874.6085
$ws.Range("H88").Value = 2394.4
$ws.Range("I88").Value = 1515.6666
$ws.Range("J88").Value = 3712.5
$ws.Range("K88").Value = 1515.6666
$ws.Range("L88").Value = 3712.5
$ws.Range("M88").Value = -1109.6666
$ws.Range("N88").Value = -4524.5
$ws.Range("H91").Value = 2394.4
$ws.Range("I91").Value = 1515.6666
$ws.Range("J91").Value = 3712.5
$ws.Range("K91").Value = 1515.6666
$ws.Range("L91").Value = 3712.5
$ws.Range("M91").Value = -111.6666
$ws.Range("N91").Value = -6520.5
$ws.Range("H116").Value = 3128.4
$ws.Range("I116").Value = 99.5
$ws.Range("K116").Value = 99.5
$ws.Range("M116").Value = 2194.5
$ws.Range("H122").Value = 2755
$ws.Range("I122").Value = 2492
$ws.Range("J122").Value = 2965.4
$ws.Range("K122").Value = 7476
$ws.Range("L122").Value = 8896.200000000001
$ws.Range("M122").Value = -5026
$ws.Range("N122").Value = -13796.2
$ws.Range("H132").Value = 4015
$ws.Range("I132").Value = 4015
$ws.Range("K132").Value = 12045
$ws.Range("M132").Value = -9515

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H3").Value = 3128.4
$ws.Range("I3").Value = 99.5
$ws.Range("K3").Value = 99.5
$ws.Range("M3").Value = 14.5
$ws.Range("H22").Value = 2216.6667
$ws.Range("I22").Value = 2259.6
$ws.Range("K22").Value = 2259.6
$ws.Range("M22").Value = -2086.6
$ws.Range("H86").Value = 1358943.2
$ws.Range("I86").Value = 1922.8889
$ws.Range("K86").Value = 1922.8889
$ws.Range("M86").Value = -799.8888999999999
$ws.Range("H89").Value = 1358943.2
$ws.Range("I89").Value = 1922.8889
$ws.Range("K89").Value = 9614.4445
$ws.Range("M89").Value = -3998.4445
$ws.Range("H99").Value = 2765.0454
$ws.Range("I99").Value = 1063.6
$ws.Range("K99").Value = 1063.6
$ws.Range("M99").Value = 434.4000000000001
$ws.Range("H132").Value = 136570.28
$ws.Range("J132").Value = 136570.28
$ws.Range("L132").Value = 136570.28
$ws.Range("N132").Value = -146690.28

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H132").Value = 1991.3846
$ws.Range("I132").Value = 1991.3846
$ws.Range("K132").Value = 5974.1538
$ws.Range("M132").Value = -3444.1538

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H4").Value = 8241615.5
$ws.Range("I4").Value = 1698200
$ws.Range("J4").Value = 78583336
$ws.Range("K4").Value = 5094600
$ws.Range("L4").Value = 235750008
$ws.Range("M4").Value = -5094488
$ws.Range("N4").Value = -235750232
$ws.Range("H55").Value = 284084.34
$ws.Range("I55").Value = 1252878.5
$ws.Range("K55").Value = 3758635.5
$ws.Range("M55").Value = -3758458.5
$ws.Range("H74").Value = 33221.332
$ws.Range("I74").Value = 26632
$ws.Range("K74").Value = 79896
$ws.Range("M74").Value = -78835
$ws.Range("H77").Value = 33221.332
$ws.Range("I77").Value = 26632
$ws.Range("K77").Value = 239688
$ws.Range("M77").Value = -234384
$ws.Range("H86").Value = 2166.6667
$ws.Range("I86").Value = 1500
$ws.Range("K86").Value = 4500
$ws.Range("M86").Value = -3314
$ws.Range("H89").Value = 2166.6667
$ws.Range("I89").Value = 1500
$ws.Range("K89").Value = 13500
$ws.Range("M89").Value = -7572

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H51").Value = 52464.465
$ws.Range("I51").Value = 44466.3
$ws.Range("J51").Value = 68460.8
$ws.Range("K51").Value = 44466.3
$ws.Range("L51").Value = 68460.8
$ws.Range("M51").Value = -43957.3
$ws.Range("N51").Value = -69478.8
$ws.Range("H122").Value = 3361.4375
$ws.Range("I122").Value = 2791.8
$ws.Range("J122").Value = 3620.3635
$ws.Range("K122").Value = 8375.400000000001
$ws.Range("L122").Value = 10861.0905
$ws.Range("M122").Value = -5925.400000000001
$ws.Range("N122").Value = -15761.0905
$ws.Range("H132").Value = 2705.1538
$ws.Range("I132").Value = 2705.1538
$ws.Range("J132").Value = 0
$ws.Range("K132").Value = 8115.4614
$ws.Range("L132").Value = 0
$ws.Range("M132").Value = -5585.4614
$ws.Range("N132").Value = $null

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H7").Value = 3799.625
$ws.Range("I7").Value = 1732.8334
$ws.Range("K7").Value = 1732.8334
$ws.Range("M7").Value = -1620.8334
$ws.Range("H40").Value = 3534.7097
$ws.Range("I40").Value = 2456.0476
$ws.Range("K40").Value = 2456.0476
$ws.Range("M40").Value = -2320.0476
$ws.Range("H55").Value = 1560.0526
$ws.Range("I55").Value = 235.90909
$ws.Range("K55").Value = 235.90909
$ws.Range("M55").Value = -62.90908999999999
$ws.Range("H122").Value = 9393
$ws.Range("I122").Value = 9366.25
$ws.Range("K122").Value = 28098.75
$ws.Range("M122").Value = -25648.75
$ws.Range("H126").Value = 3799.625
$ws.Range("I126").Value = 1732.8334
$ws.Range("K126").Value = 5198.5002
$ws.Range("M126").Value = -2728.5002
$ws.Range("I132").Value = 2720.4614
$ws.Range("J132").Value = 5250
$ws.Range("K132").Value = 8161.3842
$ws.Range("L132").Value = 15750
$ws.Range("M132").Value = -5631.3842
$ws.Range("N132").Value = -20810
$ws.Range("H133").Value = 71749.75
$ws.Range("J133").Value = 139999
$ws.Range("L133").Value = 139999
$ws.Range("N133").Value = -145059
$ws.Range("H136").Value = 25220.545
$ws.Range("I136").Value = 1941.4166
$ws.Range("J136").Value = 53155.5
$ws.Range("K136").Value = 5824.2498
$ws.Range("L136").Value = 159466.5
$ws.Range("M136").Value = -3274.2498
$ws.Range("N136").Value = -164566.5

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H96").Value = 3888.111
$ws.Range("J96").Value = 4332.5
$ws.Range("L96").Value = 4332.5
$ws.Range("N96").Value = -7078.5
$ws.Range("H122").Value = 5306.1113
$ws.Range("I122").Value = 2773.182
$ws.Range("K122").Value = 8319.545999999998
$ws.Range("M122").Value = -5869.545999999998
$ws.Range("H126").Value = 2839.8096
$ws.Range("I126").Value = 2445.4666
$ws.Range("J126").Value = 3825.6667
$ws.Range("K126").Value = 7336.399800000001
$ws.Range("L126").Value = 11477.0001
$ws.Range("M126").Value = -4866.399800000001
$ws.Range("N126").Value = -16417.0001
$ws.Range("H132").Value = 3145.9412
$ws.Range("I132").Value = 2932.8518
$ws.Range("K132").Value = 8798.555399999999
$ws.Range("M132").Value = -6268.555399999999
$ws.Range("H136").Value = 1583.6097
$ws.Range("I136").Value = 1090.0333
$ws.Range("K136").Value = 3270.0999
$ws.Range("M136").Value = -720.0999000000002
